# --- Rename sheets ---
$wb = $excel.ActiveWorkbook

$wsABC = $wb.Worksheets.Item("ABC Notes")
$wsABC.Name = "PV Notes"

$wsDEF = $wb.Worksheets.Item("DEF Notes")
$wsDEF.Name = "SP Notes"

# --- Rename tables ---
$loABC = $wsABC.ListObjects.Item(1)
$loABC.Name = "PV_NOTES"

$loDEF = $wsDEF.ListObjects.Item(1)
$loDEF.Name = "SP_NOTES"

# --- Index sheet: update SHEET_INDEX table data ---
$wsIndex = $wb.Worksheets.Item("Index")

$wsIndex.Range("J2:L6").Value = "08/2025"
$wsIndex.Range("E2:E6").Value = "TEST"
$wsIndex.Range("F2").Value = 1
$wsIndex.Range("F3").Value = 2
$wsIndex.Range("F4").Value = 3
$wsIndex.Range("F5").Value = 4
$wsIndex.Range("F6").Value = 5

# Update selection on Index sheet
$wsIndex.Range("E2").Select()

# --- DEF/SP Notes sheet: update selection ---
$wsDEF.Range("B4").Select()

# --- Excel Notes sheet: delete last row, update Sheet names ---
$wsExcelNotes = $wb.Worksheets.Item("Excel Notes")
$wsExcelNotes.Rows.Item(7).Delete()

$wsExcelNotes.Range("A2").Value = "PV-201"
$wsExcelNotes.Range("A3").Value = "PV-202"
$wsExcelNotes.Range("A4").Value = "PV-203"
$wsExcelNotes.Range("A5").Value = "PV-204"
$wsExcelNotes.Range("A6").Value = "PV-204A"

$wsExcelNotes.Range("A2").Select()

# Restore Index as the active/selected sheet (it was tabSelected in the original)
$wsIndex.Activate()

Write-Host "done"
